$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.641.65'
$ws.Range('E2').Value = '  +1.18%  '
$ws.Range('D3').Value = '1.882.02'
$ws.Range('E3').Value = '  +0.14%  '
$c = $ws.Range('D4')
$c.NumberFormat = "@"
$c.Value = '0.9988'
$c.Style = "Normal"
$ws.Range('E4').Value = '  -0.17%  '
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '239.43'
$c.Style = "Normal"
$ws.Range('E5').Value = '  +0.83%  '
$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '0.9992'
$c.Style = "Normal"
$ws.Range('E6').Value = '  -0.14%  '
$c = $ws.Range('D7')
$c.NumberFormat = "@"
$c.Value = '0.4807'
$c.Style = "Normal"
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('B8').Value = 'Cardano'
$ws.Range('C8').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$c = $ws.Range('D8')
$c.NumberFormat = "@"
$c.Value = '0.2838'
$c.Style = "Normal"
$ws.Range('E8').Value = '  -1.49%  '
$ws.Range('B9').Value = 'Dogecoin'
$ws.Range('C9').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$c = $ws.Range('D9')
$c.NumberFormat = "@"
$c.Value = '0.06550'
$c.Style = "Normal"
$ws.Range('E9').Value = '  -0.67%  '
$ws.Range('B10').Value = 'WrappedEther'
$ws.Range('C10').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D10').Value = '1.933.39'
$ws.Range('E10').Value = '  +2.59%  '
$ws.Range('B11').Value = 'TRON'
$ws.Range('C11').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$c = $ws.Range('D11')
$c.NumberFormat = "@"
$c.Value = '0.07495'
$c.Style = "Normal"
$ws.Range('E11').Value = '  +1.67%  '
$ws.Range('B12').Value = 'Solana'
$ws.Range('C12').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$c = $ws.Range('D12')
$c.NumberFormat = "@"
$c.Value = '16.71'
$c.Style = "Normal"
$ws.Range('E12').Value = '  -0.81%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$c = $ws.Range('D13')
$c.NumberFormat = "@"
$c.Value = '5.104'
$c.Style = "Normal"
$ws.Range('E13').Value = '  -1.66%  '
$ws.Range('B14').Value = 'Litecoin'
$ws.Range('C14').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$c = $ws.Range('D14')
$c.NumberFormat = "@"
$c.Value = '88.67'
$c.Style = "Normal"
$ws.Range('E14').Value = '  +0.85%  '
$ws.Range('B15').Value = 'Polygon'
$ws.Range('C15').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$c = $ws.Range('D15')
$c.NumberFormat = "@"
$c.Value = '0.6666'
$c.Style = "Normal"
$ws.Range('E15').Value = '  +0.86%  '
$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').Value = '30.597.17'
$ws.Range('E16').Value = '  +1.10%  '
$ws.Range('B17').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C17').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D17').Value = '2.416.89'
$ws.Range('E17').Value = '  +12.63%  '
$ws.Range('B18').Value = 'Avalanche'
$ws.Range('C18').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$c = $ws.Range('D18')
$c.NumberFormat = "@"
$c.Value = '13.38'
$c.Style = "Normal"
$ws.Range('E18').Value = '  -0.75%  '
$ws.Range('B19').Value = 'Dai'
$ws.Range('C19').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$c = $ws.Range('D19')
$c.NumberFormat = "@"
$c.Value = '0.9996'
$c.Style = "Normal"
$ws.Range('E19').Value = '  -0.07%  '
$ws.Range('B20').Value = 'ShibaInu'
$ws.Range('C20').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '0.000007624'
$c.Style = "Normal"
$ws.Range('E20').Value = '  -1.18%  '
$ws.Range('B21').Value = 'BitcoinCash'
$ws.Range('C21').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$c = $ws.Range('D21')
$c.NumberFormat = "@"
$c.Value = '224.62'
$c.Style = "Normal"
$ws.Range('E21').Value = '  +16.53%  '
$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '5.346'
$c.Style = "Normal"
$ws.Range('E22').Value = '  -1.99%  '
$ws.Range('B23').Value = 'BinanceUSD'
$ws.Range('C23').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '0.9989'
$c.Style = "Normal"
$ws.Range('E23').Value = '  -0.23%  '
$ws.Range('B24').Value = 'Chainlink'
$ws.Range('C24').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$c = $ws.Range('D24')
$c.NumberFormat = "@"
$c.Value = '6.244'
$c.Style = "Normal"
$ws.Range('E24').Value = '  +0.87%  '
$ws.Range('B25').Value = 'Cosmos'
$ws.Range('C25').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '9.353'
$c.Style = "Normal"
$ws.Range('E25').Value = '  -0.81%  '
$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$c = $ws.Range('D26')
$c.NumberFormat = "@"
$c.Value = '166.75'
$c.Style = "Normal"
$ws.Range('E26').Value = '  +0.39%  '
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$c = $ws.Range('D27')
$c.NumberFormat = "@"
$c.Value = '18.63'
$c.Style = "Normal"
$ws.Range('E27').Value = '  +1.70%  '
$ws.Range('B28').Value = 'LidoDAOToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$c = $ws.Range('D28')
$c.NumberFormat = "@"
$c.Value = '1.968'
$c.Style = "Normal"
$ws.Range('E28').Value = '  +1.38%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$c = $ws.Range('D29')
$c.NumberFormat = "@"
$c.Value = '1.457'
$c.Style = "Normal"
$ws.Range('E29').Value = '  +0.86%  '
$ws.Range('B30').Value = 'Stellar'
$ws.Range('C30').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value = '0.09515'
$c.Style = "Normal"
$ws.Range('E30').Value = '  +3.90%  '
$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$c = $ws.Range('D31')
$c.NumberFormat = "@"
$c.Value = '4.345'
$c.Style = "Normal"
$ws.Range('E31').Value = '  +1.85%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$c = $ws.Range('D32')
$c.NumberFormat = "@"
$c.Value = '4.043'
$c.Style = "Normal"
$ws.Range('E32').Value = '  -0.19%  '
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$c = $ws.Range('D33')
$c.NumberFormat = "@"
$c.Value = '0.05044'
$c.Style = "Normal"
$ws.Range('E33').Value = '  -0.47%  '
$ws.Range('B34').Value = 'ARBITRUM'
$ws.Range('C34').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$c = $ws.Range('D34')
$c.NumberFormat = "@"
$c.Value = '1.211'
$c.Style = "Normal"
$ws.Range('E34').Value = '  +6.06%  '
$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$c = $ws.Range('D35')
$c.NumberFormat = "@"
$c.Value = '0.7537'
$c.Style = "Normal"
$ws.Range('E35').Value = '  +1.29%  '
$ws.Range('B36').Value = 'HuobiToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$c = $ws.Range('D36')
$c.NumberFormat = "@"
$c.Value = '2.707'
$c.Style = "Normal"
$ws.Range('E36').Value = '  -0.16%  '
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$c = $ws.Range('D37')
$c.NumberFormat = "@"
$c.Value = '0.01839'
$c.Style = "Normal"
$ws.Range('E37').Value = '  +0.40%  '
$ws.Range('B38').Value = 'MXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$c = $ws.Range('D38')
$c.NumberFormat = "@"
$c.Value = '2.627'
$c.Style = "Normal"
$ws.Range('E38').Value = '  -0.22%  '
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c = $ws.Range('D39')
$c.NumberFormat = "@"
$c.Value = '2.090'
$c.Style = "Normal"
$ws.Range('E39').Value = '  +0.44%  '
$ws.Range('B40').Value = 'TrustWalletToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '0.9116'
$c.Style = "Normal"
$ws.Range('E40').Value = '  -0.35%  '
$ws.Range('B41').Value = 'Quant'
$ws.Range('C41').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$c = $ws.Range('D41')
$c.NumberFormat = "@"
$c.Value = '106.55'
$c.Style = "Normal"
$ws.Range('E41').Value = '  +0.04%  '
$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$c = $ws.Range('D42')
$c.NumberFormat = "@"
$c.Value = '5.872'
$c.Style = "Normal"
$ws.Range('E42').Value = '  -0.24%  '
$ws.Range('B43').Value = 'TheSandbox'
$ws.Range('C43').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$c = $ws.Range('D43')
$c.NumberFormat = "@"
$c.Value = '0.4299'
$c.Style = "Normal"
$ws.Range('E43').Value = '  -0.71%  '
$ws.Range('B44').Value = 'PaxDollar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$c = $ws.Range('D44')
$c.NumberFormat = "@"
$c.Value = '1.006'
$c.Style = "Normal"
$ws.Range('E44').Value = '  +0.63%  '
$ws.Range('B45').Value = 'Aptos'
$ws.Range('C45').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$c = $ws.Range('D45')
$c.NumberFormat = "@"
$c.Value = '7.501'
$c.Style = "Normal"
$ws.Range('E45').Value = '  -2.68%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$c = $ws.Range('D46')
$c.NumberFormat = "@"
$c.Value = '65.60'
$c.Style = "Normal"
$ws.Range('E46').Value = '  +0.54%  '
$ws.Range('B47').Value = 'Algorand'
$ws.Range('C47').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$c = $ws.Range('D47')
$c.NumberFormat = "@"
$c.Value = '0.1287'
$c.Style = "Normal"
$ws.Range('E47').Value = '  -5.03%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c = $ws.Range('D48')
$c.NumberFormat = "@"
$c.Value = '8.995'
$c.Style = "Normal"
$ws.Range('E48').Value = '  +0.69%  '
$ws.Range('E49').Value = '  -6.33%  '
$ws.Range('B50').Value = 'Elrond'
$ws.Range('C50').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '33.90'
$c.Style = "Normal"
$ws.Range('E50').Value = '  -0.86%  '
$ws.Range('B51').Value = 'Decentraland'
$ws.Range('C51').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '0.3913'
$c.Style = "Normal"
$ws.Range('E51').Value = '  +0.97%  '
